$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the shared-string texts referenced by column B
$ws.Range("B1").Value = "Test 1"
$ws.Range("B2").Value = "Test 2"
$ws.Range("B4").Value = "Test 4"

# Update the time values in column A
$ws.Range("A1").Value = 0.84166666666666667
$ws.Range("A2").Value = 0.84236111111111101
$ws.Range("A3").Value = 0.84375
$ws.Range("A4").Value = 0.84027777777777779
